{"js": "// Insert a new paragraph \"Thomas Novalski\" right after the existing\n// title paragraph (\"Aula Git \u2013 Arquivo Local\"), i.e. at the end of the\n// document body.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst newParagraph = lastParagraph.insertParagraph(\"Thomas Novalski\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Insert a new paragraph \"Thomas Novalski\" right after the existing\n# title paragraph (\"Aula Git \u2013 Arquivo Local\"), i.e. at the end of the\n# document body.\n$d = $word.ActiveDocument\n\n$lastParagraph = $d.Paragraphs.Last\n$lastParagraph.Range.InsertParagraphAfter()\n$newRange = $d.Paragraphs.Last.Range\n$newRange.Text = \"Thomas Novalski\"\n"}
